$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "ljfwheuodgahfdlkgjdf;kgkj"
$ws.Range("B7").Select()
